# fall 24 week 13 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("E3").Value = 10.82

$ws.Range("E4").Value = 10.61
$ws.Range("F4").Value = 9.699999999999999
$ws.Range("G4").Value = 10.05

$ws.Range("C5").Value = 9.130000000000001
$ws.Range("D5").Value = 9.390000000000001
$ws.Range("F5").Value = 10.1
$ws.Range("G5").Value = 9.720000000000001
$ws.Range("H5").Value = 8.470000000000001

$ws.Range("D6").Value = 10.3
$ws.Range("E6").Value = 9.9
$ws.Range("G6").Value = 10.33
$ws.Range("H6").Value = 10.55

$ws.Range("D7").Value = 10.11
$ws.Range("E7").Value = 10.28
$ws.Range("F7").Value = 9.67
$ws.Range("H7").Value = 9.970000000000001
$ws.Range("J7").Value = 9.779999999999999

$ws.Range("E8").Value = 11.53
$ws.Range("F8").Value = 9.449999999999999
$ws.Range("G8").Value = 10.03

$ws.Range("G10").Value = 10.22
